$d = $word.ActiveDocument

# Delete paragraphs 2 through 7 (the "rd commit" paragraph through "da"),
# including their paragraph marks, leaving paragraph 1 ("3") and the
# last paragraph ("sd") intact.
$start = $d.Paragraphs(2).Range.Start
$end = $d.Paragraphs(8).Range.Start
$r = $d.Range($start, $end)
$r.Delete()

# Move the "_GoBack" bookmark from the end of the final paragraph ("sd")
# to its start.
$last = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks("_GoBack").Delete()
$pos = $last.Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos))
